$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" = strikeouts) is being regenerated from "Strike#" to "K" stat.
# Update the K column values for each row (G5 unchanged).
$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 0
